$d = $word.ActiveDocument

# Locate the target paragraph: "To build another new type of NN: The Deep
# Tree-Based RNN that makes use of a tree based NN to understand the
# contents of python programs."
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Deep Tree-Based RNN*") {
        $target = $p
    }
}

$pStart = $target.Range.Start
$pText = $target.Range.Text

# The substring spanning the four runs we need to fold into one: " " +
# "that makes use of a " + "tree" + " based NN to understand the contents
# of python programs". Replacing it with its own (identical) text is a
# no-op textually, but it makes the engine rebuild the run list for the
# matched span as a single run sharing the first matched run's
# formatting.
$find = " that makes use of a tree based NN to understand the contents of python programs"
$offset = $pText.IndexOf($find)
$midStart = $pStart + $offset
$dotStart = $midStart + $find.Length

$d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $find, 2)

# That replace merges the whole rest of the paragraph (including the
# trailing ".") into one run because every run in the paragraph shares
# identical rPr. Re-split the trailing "." back into its own run, and
# re-split the leading "To build ... RNN" run away from the merged
# middle run, by toggling a formatting property on/off on each piece -
# this forces the engine to split the run at that boundary without
# altering the effective formatting (the toggle is reverted before the
# document is saved).
$dotEnd = $dotStart + 1
$rDot = $d.Range($dotStart, $dotEnd)
$rDot.Bold = 1
$rDot.Bold = 0

$rMid = $d.Range($midStart, $dotStart)
$rMid.Bold = 1
$rMid.Bold = 0
